# "Day 2 Lab Challenge / Build List.xlsx" re-edit
#
# Summary of the change being applied (per the supplied diff):
#   - Details sheet: the "Gateway" row (10.0.1.0/24) is removed from the
#     Virtual-Network table (Table1), shrinking it from B4:D8 to B4:D7 and
#     shifting the Virtual-Machine table (Table2) and the internal-load-
#     balancer row up by one row.
#   - Details sheet: the whole "Storage" table (Table3, B18:D20 originally)
#     is deleted.
#   - The Details sheet becomes the active/selected sheet (instead of NSG),
#     with a new zoom level and a new active cell.
#
# All of the above is performed with native Excel object-model calls so the
# engine recomputes sharedStrings.xml, table refs/ids, dimensions, etc. on
# save exactly the way Excel would after these UI actions.

$wb = $excel.ActiveWorkbook
$details = $wb.Worksheets.Item("Details")

# --- Remove the "Gateway" row from the Virtual Network table -------------
# Row 7 holds: C7 = "Gateway", D7 = "10.0.1.0/24". Deleting the whole sheet
# row shifts every row below it (the LAN row, and both the Virtual Machine
# table and the Storage table) up by one, exactly like the diff shows.
$details.Rows.Item(7).Delete()

# --- Remove the Storage table entirely ------------------------------------
# After the shift above, the Storage table (Storage/Type/Purpose) now sits
# at B17:D19. Deleting the ListObject clears its cells; Excel then drops the
# now-empty trailing rows from the sheet's used range automatically.
$details.ListObjects.Item("Table3").Delete()

# --- Make "Details" the active sheet with the new selection/zoom ---------
$details.Activate()
$excel.ActiveWindow.Zoom = 160
$details.Range("G10").Select()
